{"js": "// Update the date line and the 25 \"three-digit \u00d7 one-digit\" practice\n// answers in the single table, in document order, leaving all other\n// formatting untouched.\n\n// 1) Update the date paragraph above the table.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.insertText(\"2024-02-29 Thursday\", \"Replace\");\n\n// 2) Update the equation table. The table has 20 rows (5 data rows with\n// 4 blank spacer rows interleaved); data lives in rows 0, 4, 9, 14, 19.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// New values in document order (row-major across the 5 data rows).\nconst newValues = [\n  [\"143\u00d79=1287\", \"399\u00d77=2793\", \"142\u00d75=710\", \"805\u00d78=6440\", \"885\u00d72=1770\"],\n  [\"588\u00d76=3528\", \"183\u00d73=549\", \"414\u00d79=3726\", \"922\u00d74=3688\", \"227\u00d75=1135\"],\n  [\"118\u00d75=590\", \"847\u00d79=7623\", \"953\u00d77=6671\", \"218\u00d75=1090\", \"531\u00d72=1062\"],\n  [\"304\u00d75=1520\", \"866\u00d78=6928\", \"518\u00d79=4662\", \"976\u00d72=1952\", \"638\u00d73=1914\"],\n  [\"128\u00d78=1024\", \"929\u00d76=5574\", \"101\u00d79=909\", \"807\u00d72=1614\", \"325\u00d79=2925\"],\n];\n\nconst dataRowIndexes = [0, 4, 9, 14, 19];\n\nfor (let i = 0; i < dataRowIndexes.length; i++) {\n  const rowIndex = dataRowIndexes[i];\n  for (let col = 0; col < newValues[i].length; col++) {\n    table.getCell(rowIndex, col).value = newValues[i][col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 \"three-digit x one-digit\" practice\n# answers in the single table, in document order, leaving all other\n# formatting untouched.\n\n$d = $word.ActiveDocument\n\n# 1) Update the date paragraph above the table.\n$d.Paragraphs(1).Range.Text = \"2024-02-29 Thursday\"\n\n# 2) Update the equation table. The table has 20 rows (5 data rows with\n# 4 blank spacer rows interleaved); data lives in rows 1, 5, 10, 15, 20\n# (1-based Word COM indexing).\n$t = $d.Tables(1)\n\n$newValues = @(\n    @(\"143\u00d79=1287\", \"399\u00d77=2793\", \"142\u00d75=710\", \"805\u00d78=6440\", \"885\u00d72=1770\"),\n    @(\"588\u00d76=3528\", \"183\u00d73=549\", \"414\u00d79=3726\", \"922\u00d74=3688\", \"227\u00d75=1135\"),\n    @(\"118\u00d75=590\", \"847\u00d79=7623\", \"953\u00d77=6671\", \"218\u00d75=1090\", \"531\u00d72=1062\"),\n    @(\"304\u00d75=1520\", \"866\u00d78=6928\", \"518\u00d79=4662\", \"976\u00d72=1952\", \"638\u00d73=1914\"),\n    @(\"128\u00d78=1024\", \"929\u00d76=5574\", \"101\u00d79=909\", \"807\u00d72=1614\", \"325\u00d79=2925\")\n)\n\n$dataRowIndexes = @(1, 5, 10, 15, 20)\n\nfor ($i = 0; $i -lt $dataRowIndexes.Count; $i++) {\n    $rowIndex = $dataRowIndexes[$i]\n    for ($col = 1; $col -le 5; $col++) {\n        $t.Cell($rowIndex, $col).Range.Text = $newValues[$i][$col - 1]\n    }\n}\n"}
